# 7.62x54r FMJ and 7.92 FMJ price and perf matching
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# ammo_7.92x33_fmj (row 4): update price
$ws.Range("C4").Value = 2250

# ammo_7.62x54_7h1 (row 6): update price and pen value
$ws.Range("C6").Value = 2250
$ws.Range("G6").Value = 0.3

# ammo_7.62x54_ap (row 7): update price
$ws.Range("C7").Value = 9500

# Update selection to match the saved cursor position in the workbook
$ws.Activate()
$ws.Range("D6").Select()
